# Rename sheets and update the task-order file references for the
# "Practice tasks and final revisions" commit.

$wb = $excel.ActiveWorkbook

# --- Sheet renames ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16509960819059103"
$wb.Worksheets.Item(2).Name = "NB_TO-16509960838507922"
$wb.Worksheets.Item(3).Name = "RS_TO-16509960838507922"
$wb.Worksheets.Item(4).Name = "TOL_TO-16509960838988266"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1650996083962824"

# --- Sheet 1 (GNG) updated stim file names ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1650996081873912.csv"
$ws1.Range("B3").Value = "GNG_stims-16509960818899698.csv"
$ws1.Range("B4").Value = "go_stims-16509960818899698.csv"
$ws1.Range("B5").Value = "GNG_stims-16509960819059103.csv"

# --- Sheet 2 (NB) updated stim file names ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16509960835947926.csv"
$ws2.Range("B3").Value = "ZB-match_8-16509960822507892.csv"
$ws2.Range("B4").Value = "TB-16509960837308228.csv"
$ws2.Range("B5").Value = "ZB-match_7-16509960820988262.csv"
$ws2.Range("B6").Value = "OB-16509960828187945.csv"
$ws2.Range("B7").Value = "OB-16509960826028001.csv"
$ws2.Range("B8").Value = "ZB-match_5-16509960825468018.csv"
$ws2.Range("B9").Value = "TB-16509960838268182.csv"
$ws2.Range("B10").Value = "OB-16509960834588234.csv"

# --- Sheet 3 (RS) has no data changes, only the rename above ---

# --- Sheet 4 (TOL) updated stim file names ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-1650996083866796.csv"
$ws4.Range("B3").Value = "ZM_stims-16509960838507922.csv"
$ws4.Range("B4").Value = "MM_stims-16509960838827925.csv"
$ws4.Range("B5").Value = "ZM_stims-1650996083866796.csv"
$ws4.Range("B6").Value = "MM_stims-16509960838988266.csv"
$ws4.Range("B7").Value = "ZM_stims-16509960838827925.csv"

# --- Sheet 5 (vSAT) updated stim file names ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16509960838988266.csv"
$ws5.Range("B3").Value = "vSAT_stims-16509960839308276.csv"
$ws5.Range("B4").Value = "SAT_stims-1650996083914826.csv"
$ws5.Range("B5").Value = "vSAT_stims-16509960839468222.csv"
